$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(249, 1).Value = "2024-12-19 00:30:01"
$ws.Cells.Item(249, 2).Value = -0.119228246283019
$ws.Cells.Item(249, 3).Value = -0.001977707633431996
$ws.Cells.Item(249, 4).Value = 0.009431944511785462

$ws.Cells.Item(250, 1).Value = "2024-12-19 00:30:02"
$ws.Cells.Item(250, 2).Value = -0.1203799240972148
$ws.Cells.Item(250, 3).Value = -0.002208202723795996
$ws.Cells.Item(250, 4).Value = 0.010632931051273

$ws.Cells.Item(251, 1).Value = "2024-12-19 00:30:03"
$ws.Cells.Item(251, 2).Value = -0.1205848837082158
$ws.Cells.Item(251, 3).Value = -0.002008025279265996
$ws.Cells.Item(251, 4).Value = 0.009685499791337909

$ws.Cells.Item(252, 1).Value = "2024-12-19 00:30:04"
$ws.Cells.Item(252, 2).Value = -0.119712992029672
$ws.Cells.Item(252, 3).Value = -0.002061119119799996
$ws.Cells.Item(252, 4).Value = 0.009869709470432862

$ws.Cells.Item(253, 1).Value = "2024-12-19 00:30:05"
$ws.Cells.Item(253, 2).Value = -0.1203506441527861
$ws.Cells.Item(253, 3).Value = -0.002076708159727995
$ws.Cells.Item(253, 4).Value = 0.009997326589624451

$ws.Cells.Item(254, 1).Value = "2024-12-19 00:30:06"
$ws.Cells.Item(254, 2).Value = -0.1196153922149097
$ws.Cells.Item(254, 3).Value = -0.002024120456853995
$ws.Cells.Item(254, 4).Value = 0.00968463849347251

$ws.Cells.Item(255, 1).Value = "2024-12-19 00:30:07"
$ws.Cells.Item(255, 2).Value = -0.1193616326965275
$ws.Cells.Item(255, 3).Value = -0.001980592618093995
$ws.Cells.Item(255, 4).Value = 0.009456270744095569

$ws.Cells.Item(256, 1).Value = "2024-12-19 00:30:08"
$ws.Cells.Item(256, 2).Value = -0.1202660576466587
$ws.Cells.Item(256, 3).Value = -0.001892727120317995
$ws.Cells.Item(256, 4).Value = 0.009105233158462338

$ws.Cells.Item(257, 1).Value = "2024-12-19 00:30:09"
$ws.Cells.Item(257, 2).Value = -0.119752031955577
$ws.Cells.Item(257, 3).Value = -0.002054589943985996
$ws.Cells.Item(257, 4).Value = 0.009841652825112724

$ws.Cells.Item(258, 1).Value = "2024-12-19 00:30:10"
$ws.Cells.Item(258, 2).Value = -0.1192054729929077
$ws.Cells.Item(258, 3).Value = -0.002058942727861995
$ws.Cells.Item(258, 4).Value = 0.009817489669603876

$ws.Cells.Item(259, 1).Value = "2024-12-19 00:30:11"
$ws.Cells.Item(259, 2).Value = -0.118691447301826
$ws.Cells.Item(259, 3).Value = -0.002134660921797996
$ws.Cells.Item(259, 4).Value = 0.01013463977227416

$ws.Cells.Item(260, 1).Value = "2024-12-19 00:30:12"
$ws.Cells.Item(260, 2).Value = -0.118307554697094
$ws.Cells.Item(260, 3).Value = -0.001991980715443996
$ws.Cells.Item(260, 4).Value = 0.009426654697917881

$ws.Cells.Item(261, 1).Value = "2024-12-19 00:30:13"
$ws.Cells.Item(261, 2).Value = -0.1194885124557186
$ws.Cells.Item(261, 3).Value = -0.002062890601609996
$ws.Cells.Item(261, 4).Value = 0.009859669173810433

$ws.Cells.Item(262, 1).Value = "2024-12-19 00:30:14"
$ws.Cells.Item(262, 2).Value = -0.1213298956275685
$ws.Cells.Item(262, 3).Value = -0.001953868549645996
$ws.Cells.Item(262, 4).Value = 0.009482506687941495

$ws.Cells.Item(263, 1).Value = "2024-12-19 00:30:15"
$ws.Cells.Item(263, 2).Value = -0.1202660576466587
$ws.Cells.Item(263, 3).Value = -0.001996789023213996
$ws.Cells.Item(263, 4).Value = 0.009605837750962793

$ws.Cells.Item(264, 1).Value = "2024-12-19 00:30:16"
$ws.Cells.Item(264, 2).Value = -0.1196967253938783
$ws.Cells.Item(264, 3).Value = -0.002068458115869996
$ws.Cells.Item(264, 4).Value = 0.009903506523361192

$ws.Cells.Item(265, 1).Value = "2024-12-19 00:30:17"
$ws.Cells.Item(265, 2).Value = -0.1200741113442928
$ws.Cells.Item(265, 3).Value = -0.001968900838147995
$ws.Cells.Item(265, 4).Value = 0.009456560738626149

$ws.Cells.Item(266, 1).Value = "2024-12-19 00:30:18"
$ws.Cells.Item(266, 2).Value = -0.1197487786284182
$ws.Cells.Item(266, 3).Value = -0.001942379224763994
$ws.Cells.Item(266, 4).Value = 0.009303901591948088

$ws.Cells.Item(267, 1).Value = "2024-12-19 00:30:19"
$ws.Cells.Item(267, 2).Value = -0.1197878185543232
$ws.Cells.Item(267, 3).Value = -0.002152527581195997
$ws.Cells.Item(267, 4).Value = 0.01031386333317929

$ws.Cells.Item(268, 1).Value = "2024-12-19 00:30:20"
$ws.Cells.Item(268, 2).Value = -0.1208646698438679
$ws.Cells.Item(268, 3).Value = -0.001925119930557995
$ws.Cells.Item(268, 4).Value = 0.009307159392669679

$ws.Cells.Item(269, 1).Value = "2024-12-19 00:30:21"
$ws.Cells.Item(269, 2).Value = -0.1211737359239487
$ws.Cells.Item(269, 3).Value = -0.001908569229075996
$ws.Cells.Item(269, 4).Value = 0.009250738550265165

$ws.Cells.Item(270, 1).Value = "2024-12-19 00:30:22"
$ws.Cells.Item(270, 2).Value = -0.1220716542197624
$ws.Cells.Item(270, 3).Value = -0.001968698383083996
$ws.Cells.Item(270, 4).Value = 0.009612890731313396

$ws.Cells.Item(271, 1).Value = "2024-12-19 00:30:23"
$ws.Cells.Item(271, 2).Value = -0.1215771484916331
$ws.Cells.Item(271, 3).Value = -0.001928915963007995
$ws.Cells.Item(271, 4).Value = 0.009380484098500182

$ws.Cells.Item(272, 1).Value = "2024-12-19 00:30:24"
$ws.Cells.Item(272, 2).Value = -0.121681254960713
$ws.Cells.Item(272, 3).Value = -0.001945871574617996
$ws.Cells.Item(272, 4).Value = 0.009471043807675853

$ws.Cells.Item(273, 1).Value = "2024-12-19 00:30:25"
$ws.Cells.Item(273, 2).Value = -0.1226019465466379
$ws.Cells.Item(273, 3).Value = -0.001862207019419996
$ws.Cells.Item(273, 4).Value = 0.009132408218148171

$ws.Cells.Item(274, 1).Value = "2024-12-19 00:30:26"
$ws.Cells.Item(274, 2).Value = -0.1199537382394192
$ws.Cells.Item(274, 3).Value = -0.002046896651553995
$ws.Cells.Item(274, 4).Value = 0.009821316205746061

$ws.Cells.Item(275, 1).Value = "2024-12-19 00:30:27"
$ws.Cells.Item(275, 2).Value = -0.1219935743679525
$ws.Cells.Item(275, 3).Value = -0.002187552307267995
$ws.Cells.Item(275, 4).Value = 0.01067469300321937

$ws.Cells.Item(276, 1).Value = "2024-12-19 00:30:28"
$ws.Cells.Item(276, 2).Value = -0.1206044036711683
$ws.Cells.Item(276, 3).Value = -0.002062131395119996
$ws.Cells.Item(276, 4).Value = 0.009948085088001657

$ws.Cells.Item(277, 1).Value = "2024-12-19 00:30:29"
$ws.Cells.Item(277, 2).Value = -0.121310375664616
$ws.Cells.Item(277, 3).Value = -0.001907455726223996
$ws.Cells.Item(277, 4).Value = 0.009255766828474234

$ws.Cells.Item(278, 1).Value = "2024-12-19 00:30:30"
$ws.Cells.Item(278, 2).Value = -0.1200903779800865
$ws.Cells.Item(278, 3).Value = -0.002034597506415997
$ws.Cells.Item(278, 4).Value = 0.009773423343313541

$ws.Cells.Item(279, 1).Value = "2024-12-19 00:30:31"
$ws.Cells.Item(279, 2).Value = -0.1205978970168508
$ws.Cells.Item(279, 3).Value = -0.002062080781353996
$ws.Cells.Item(279, 4).Value = 0.009947304228406255

$ws.Cells.Item(280, 1).Value = "2024-12-19 00:30:32"
$ws.Cells.Item(280, 2).Value = -0.1207573100476293
$ws.Cells.Item(280, 3).Value = -0.001983123306393996
$ws.Cells.Item(280, 4).Value = 0.009579065438915979

$ws.Cells.Item(281, 1).Value = "2024-12-19 00:30:33"
$ws.Cells.Item(281, 2).Value = -0.1203994440601673
$ws.Cells.Item(281, 3).Value = -0.001964345599207996
$ws.Cells.Item(281, 4).Value = 0.009460244723467156

$ws.Cells.Item(282, 1).Value = "2024-12-19 00:30:34"
$ws.Cells.Item(282, 2).Value = -0.1197943252086407
$ws.Cells.Item(282, 3).Value = -0.001825512039069996
$ws.Cells.Item(282, 4).Value = 0.008747439315225593

$ws.Cells.Item(283, 1).Value = "2024-12-19 00:30:35"
$ws.Cells.Item(283, 2).Value = -0.1205523504366283
$ws.Cells.Item(283, 3).Value = -0.002044922714679996
$ws.Cells.Item(283, 4).Value = 0.009860809588636969

$ws.Cells.Item(284, 1).Value = "2024-12-19 00:30:36"
$ws.Cells.Item(284, 2).Value = -0.1204937905477709
$ws.Cells.Item(284, 3).Value = -0.001753589877583996
$ws.Cells.Item(284, 4).Value = 0.008451867656651889

$ws.Cells.Item(285, 1).Value = "2024-12-19 00:30:37"
$ws.Cells.Item(285, 2).Value = -0.1194397125483374
$ws.Cells.Item(285, 3).Value = -0.002023513091661996
$ws.Cells.Item(285, 4).Value = 0.009667512880236253

$ws.Cells.Item(286, 1).Value = "2024-12-19 00:30:38"
$ws.Cells.Item(286, 2).Value = -0.1209980562573764
$ws.Cells.Item(286, 3).Value = -0.001985755222225996
$ws.Cells.Item(286, 4).Value = 0.009610900883691207

$ws.Cells.Item(287, 1).Value = "2024-12-19 00:30:39"
$ws.Cells.Item(287, 2).Value = -0.1216519750162843
$ws.Cells.Item(287, 3).Value = -0.002059094569159996
$ws.Cells.Item(287, 4).Value = 0.01001971684334474

$ws.Cells.Item(288, 1).Value = "2024-12-19 00:30:40"
$ws.Cells.Item(288, 2).Value = -0.1190720865793992
$ws.Cells.Item(288, 3).Value = -0.002046491741425996
$ws.Cells.Item(288, 4).Value = 0.009747201672764065

$ws.Cells.Item(289, 1).Value = "2024-12-19 00:30:41"
$ws.Cells.Item(289, 2).Value = -0.1177187024813611
$ws.Cells.Item(289, 3).Value = -0.001913478764377996
$ws.Cells.Item(289, 4).Value = 0.00901008949472863

$ws.Cells.Item(290, 1).Value = "2024-12-19 00:30:42"
$ws.Cells.Item(290, 2).Value = -0.1179301687466796
$ws.Cells.Item(290, 3).Value = -0.001976594130579996
$ws.Cells.Item(290, 4).Value = 0.009324003174519814

$ws.Cells.Item(291, 1).Value = "2024-12-19 00:30:43"
$ws.Cells.Item(291, 2).Value = -0.1194266992397024
$ws.Cells.Item(291, 3).Value = -0.002047453402979996
$ws.Cells.Item(291, 4).Value = 0.009780824070599888

$ws.Cells.Item(292, 1).Value = "2024-12-19 00:30:44"
$ws.Cells.Item(292, 2).Value = -0.1194039259495912
$ws.Cells.Item(292, 3).Value = -0.002022652657639996
$ws.Cells.Item(292, 4).Value = 0.009660506726183595

$ws.Cells.Item(293, 1).Value = "2024-12-19 00:30:45"
$ws.Cells.Item(293, 2).Value = -0.1196576854679733
$ws.Cells.Item(293, 3).Value = -0.002090070193951997
$ws.Cells.Item(293, 4).Value = 0.01000371847495576

$ws.Cells.Item(294, 1).Value = "2024-12-19 00:30:46"
$ws.Cells.Item(294, 2).Value = -0.1185352875982062
$ws.Cells.Item(294, 3).Value = -0.001958626243649996
$ws.Cells.Item(294, 4).Value = 0.009286653003537863

$ws.Cells.Item(295, 1).Value = "2024-12-19 00:30:47"
$ws.Cells.Item(295, 2).Value = -0.1176601425925037
$ws.Cells.Item(295, 3).Value = -0.001845048952745997
$ws.Cells.Item(295, 4).Value = 0.008683548914809743
